# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.370.08"
$ws.Range("E2").Value = "  -2.43%  "

$ws.Range("D3").Value = "3.556.17"
$ws.Range("E3").Value = "  -2.67%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'583.91"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "'180.57"
$ws.Range("E6").Value = "  +1.86%  "

$ws.Range("D7").Value = "3.550.38"
$ws.Range("E7").Value = "  -2.77%  "

$ws.Range("E8").Value = "  -3.59%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").Value = "'0.666"
$ws.Range("E10").Value = "  -5.96%  "

$ws.Range("E11").Value = "  -10.58%  "

$ws.Range("D12").Value = "'53.34"
$ws.Range("E12").Value = "  -3.02%  "

$ws.Range("E13").Value = "  -13.65%  "

$ws.Range("D14").Value = "'9.83"
$ws.Range("E14").Value = "  -7.17%  "

$ws.Range("D15").Value = "4.140.75"
$ws.Range("E15").Value = "  -2.12%  "

$ws.Range("D16").Value = "3.560.52"
$ws.Range("E16").Value = "  -2.57%  "

$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.16"
$ws.Range("E18").Value = "  -5.43%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "66.156.45"
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("D20").Value = "'12.06"
$ws.Range("E20").Value = "  -4.86%  "

$ws.Range("E21").Value = "  -6.31%  "

$ws.Range("D22").Value = "'389.87"
$ws.Range("E22").Value = "  -4.38%  "

$ws.Range("E23").Value = "  -6.23%  "

$ws.Range("D24").Value = "'84.30"
$ws.Range("E24").Value = "  -4.50%  "

$ws.Range("D25").Value = "'2.85"
$ws.Range("E25").Value = "  -4.42%  "

$ws.Range("D26").Value = "'12.13"
$ws.Range("E26").Value = "  -3.42%  "

$ws.Range("D27").Value = "'6.03"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  -5.47%  "

$ws.Range("E29").Value = "  -7.45%  "

$ws.Range("D30").Value = "'8.87"
$ws.Range("E30").Value = "  -6.60%  "

$ws.Range("D31").Value = "'30.86"
$ws.Range("E31").Value = "  -5.03%  "

$ws.Range("E32").Value = "  -7.34%  "

$ws.Range("D33").Value = "'65.08"
$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").Value = "'11.84"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'593.66"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("E36").Value = "  -4.57%  "

$ws.Range("D37").Value = "'41.20"
$ws.Range("E37").Value = "  -3.26%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("E40").Value = "  -6.55%  "

$ws.Range("D41").Value = "0.0₃0729"
$ws.Range("E41").Value = "  -17.20%  "

$ws.Range("E42").Value = "  -5.26%  "

$ws.Range("E43").Value = "  -8.66%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.779.04"
$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0406"
$ws.Range("E45").Value = "  -6.93%  "

$ws.Range("D46").Value = "'2.39"
$ws.Range("E46").Value = "  -10.44%  "

$ws.Range("D47").Value = "'3.09"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  -3.43%  "

$ws.Range("E49").Value = "  -7.15%  "

$ws.Range("D50").Value = "'134.80"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("D51").Value = "'8.15"
$ws.Range("E51").Value = "  -8.55%  "

